$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-24 Saturday" "2024-02-25 Sunday"

Replace-Text "62÷2=31, 0" "88÷4=22, 0"
Replace-Text "30÷5=6, 0" "94÷9=10, 4"
Replace-Text "33÷5=6, 3" "48÷6=8, 0"
Replace-Text "80÷3=26, 2" "64÷5=12, 4"
Replace-Text "66÷8=8, 2" "43÷4=10, 3"
Replace-Text "81÷8=10, 1" "63÷4=15, 3"
Replace-Text "23÷2=11, 1" "74÷6=12, 2"
Replace-Text "89÷3=29, 2" "87÷7=12, 3"
Replace-Text "72÷8=9, 0" "66÷7=9, 3"
Replace-Text "48÷5=9, 3" "98÷2=49, 0"
Replace-Text "94÷5=18, 4" "41÷6=6, 5"
Replace-Text "12÷8=1, 4" "53÷3=17, 2"
Replace-Text "91÷3=30, 1" "31÷5=6, 1"
Replace-Text "90÷2=45, 0" "87÷5=17, 2"
Replace-Text "71÷9=7, 8" "50÷9=5, 5"
Replace-Text "16÷8=2, 0" "40÷5=8, 0"
Replace-Text "52÷4=13, 0" "27÷6=4, 3"
Replace-Text "10÷9=1, 1" "34÷7=4, 6"
Replace-Text "37÷8=4, 5" "43÷3=14, 1"
Replace-Text "17÷2=8, 1" "68÷2=34, 0"
Replace-Text "26÷9=2, 8" "87÷8=10, 7"
Replace-Text "45÷7=6, 3" "71÷7=10, 1"
Replace-Text "20÷2=10, 0" "12÷6=2, 0"
Replace-Text "32÷3=10, 2" "35÷5=7, 0"
Replace-Text "52÷2=26, 0" "76÷4=19, 0"

Write-Output "Done"
